# Making grid compatible with Edge
# - Add "Edge" as a new option alongside "FireFox" in the Browser column.
# - TestCases row 2 (Login_02 / Login into the site): RunMode -> Yes, Browser -> Edge
# - TestCases row 4 (Login_03 / Login into the site): RunMode -> No
# - TestCases sheet becomes the active/selected sheet & selected cell becomes D2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Update the data grid values
$ws.Range("C2").Value = "Yes"
$ws.Range("D2").Value = "Edge"
$ws.Range("C4").Value = "No"

# Make TestCases the active sheet with D2 selected (matches the authored view state)
$ws.Activate() | Out-Null
$ws.Range("D2").Select() | Out-Null
